$d = $word.ActiveDocument
$t = $d.Tables(1)
$values = @(
  "3+2=5",
  "48-8=40",
  "79+6=85",
  "60+24=84",
  "77-76=1",
  "28+17=45",
  "25+50=75",
  "56-56=0",
  "32+38=70",
  "28+10=38",
  "78+11=89",
  "67-54=13",
  "9+72=81",
  "62-48=14",
  "2+5=7",
  "72-42=30",
  "38+13=51",
  "90+6=96",
  "14+49=63",
  "94-59=35",
  "71-8=63",
  "22+30=52",
  "27+69=96",
  "9-8=1",
  "52-40=12",
  "32-9=23",
  "85-14=71",
  "48+36=84",
  "99-46=53",
  "14+38=52",
  "62+7=69",
  "55-24=31",
  "32-5=27",
  "87-6=81",
  "53-2=51",
  "64-16=48",
  "1+71=72",
  "37+55=92",
  "34-20=14",
  "32+58=90",
  "36+8=44",
  "81-71=10",
  "77+4=81",
  "14+11=25",
  "86-4=82",
  "49-39=10",
  "61+26=87",
  "73+19=92",
  "5+23=28",
  "48-35=13",
  "82-12=70",
  "57-47=10",
  "37-31=6",
  "18+3=21",
  "62+32=94",
  "23-17=6",
  "66-56=10",
  "95-8=87",
  "80-11=69",
  "84-34=50",
  "6+88=94",
  "92-87=5",
  "90-72=18",
  "42-15=27",
  "18+38=56",
  "80-73=7",
  "55-11=44",
  "24+58=82",
  "50-35=15",
  "38-18=20",
  "88-34=54",
  "10+59=69",
  "62-23=39",
  "55-42=13",
  "36+59=95",
  "38+12=50",
  "46-31=15",
  "30-5=25",
  "40+24=64",
  "56+13=69",
  "3+56=59",
  "96-81=15",
  "5+75=80",
  "43-9=34",
  "94-15=79",
  "0+37=37",
  "42-38=4",
  "23+14=37",
  "94+2=96",
  "22+40=62",
  "57+22=79",
  "67+2=69",
  "68+13=81",
  "71-55=16",
  "11+28=39",
  "2+90=92",
  "62-28=34",
  "50+38=88",
  "66+2=68",
  "44+1=45"
)
$i = 0
for ($r = 1; $r -le $t.Rows.Count; $r++) {
  for ($c = 1; $c -le $t.Columns.Count; $c++) {
    $cell = $t.Cell($r, $c)
    $cell.Range.Text = $values[$i]
    $i++
  }
}
Write-Output "updated $i cells"
